# Generate Report for Handoff
# Stamp a fresh "Latest Handoff Datetime" onto the most recently handed-off
# file (row 5, the de2e3f76-... file) on each locale sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-22 16:35:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-22 16:35:29"
